# Weekly update: prepend the latest week's two "Fruta, Femacal de La Calera - Limón"
# price records (1a amarillo / 2a amarillo) ahead of the existing historical rows,
# pushing all prior data rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 1099 (first row of the
# block that needs to move down). Everything from old row 1099 onward shifts to
# row+2.
$ws.Rows("1099:1100").Insert()

# The rows that used to be 1099 and 1100 are now at 1101 and 1102. Clone their
# full contents into the freshly-inserted blank rows 1099/1100 so every column
# (Mercado ID, Mercado, Región, Codreg, Tipo, Producto, Categoría, Variedad,
# Calidad, Unidad, Origen, Kg/unidad, ...) matches, then overwrite just the
# columns that actually hold new-week data (Fecha, Volumen, Precio mínimo,
# Precio máximo, Precio promedio ponderado, Precio $/Kg).
$ws.Range("A1101:T1101").Copy()
$ws.Range("A1099:T1099").PasteSpecial()

$ws.Range("A1102:T1102").Copy()
$ws.Range("A1100:T1100").PasteSpecial()

$excel.CutCopyMode = 0

# New week's figures for the "1a amarillo" quality row (now row 1099).
$ws.Range("D1099").Value = 44826
$ws.Range("M1099").Value = 312
$ws.Range("N1099").Value = 3500
$ws.Range("O1099").Value = 4000
$ws.Range("P1099").Value = 3768
$ws.Range("S1099").Value = 236

# New week's figures for the "2a amarillo" quality row (now row 1100).
$ws.Range("D1100").Value = 44826
$ws.Range("M1100").Value = 320
$ws.Range("N1100").Value = 2500
$ws.Range("O1100").Value = 3000
$ws.Range("P1100").Value = 2781
$ws.Range("S1100").Value = 174
